# Apply cryptocurrency price/volume updates scraped on Sun Feb 18 23:08:47 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.223.19"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "2.882.56"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.558"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0854"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").Value = "3.333.93"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.997"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.52%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.873.15"
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("D18").Value = "52.226.58"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.01%  "
$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0901"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.38%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0454"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("E40").Value = "  +3.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.85%  "
$ws.Range("D47").Value = "2.177.73"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("E48").Value = "  +6.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.242"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +14.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.965"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.85%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.07%  "
